# Update automatico via Actualizar 03-08-2021 13-23-19
# Shifts the "Ultimo" timestamp history down one block and stamps the
# newest check time on rows 2-15, matching the rolling-log layout of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = 44263.55762613597
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = 44263.53626993056
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = 44263.51489300926
}
